$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.681.92"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "2.467.23"
$ws.Range("E3").Value = "  -0.76%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'316.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.95%  "
$ws.Range("D6").Value = "'92.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("E7").Value = "  +0.57%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "'0.514"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.21%  "
$ws.Range("D10").Value = "'32.80"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.59%  "
$ws.Range("D11").Value = "'0.0843"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +7.46%  "
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("D13").Value = "2.848.34"
$ws.Range("E13").Value = "  -0.62%  "
$ws.Range("D14").Value = "'6.89"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.79%  "
$ws.Range("D15").Value = "'15.74"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.06%  "
$ws.Range("D16").Value = "2.465.16"
$ws.Range("E16").Value = "  -1.28%  "
$ws.Range("D17").Value = "'0.783"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.85%  "
$ws.Range("D18").Value = "41.665.44"
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("D19").Value = "'6.50"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.38%  "
$ws.Range("D20").Value = "0.0₃0951"
$ws.Range("E20").Value = "  +2.16%  "
$ws.Range("B21").Value = "Litecoin"
$ws.Range("C21").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D21").Value = "'71.12"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.46%  "
$ws.Range("B22").Value = "InternetComputer(DFINITY)"
$ws.Range("C22").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D22").Value = "'11.48"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.39%  "
$ws.Range("D23").Value = "'239.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.39%  "
$ws.Range("E24").Value = "  +0.95%  "
$ws.Range("E25").Value = "  +1.38%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").Value = "'24.72"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.21%  "
$ws.Range("E28").Value = "  +1.95%  "
$ws.Range("D29").Value = "'9.81"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.59%  "
$ws.Range("D30").Value = "'35.57"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.13%  "
$ws.Range("D31").Value = "'156.10"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.60%  "
$ws.Range("D32").Value = "'5.50"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.11%  "
$ws.Range("E33").Value = "  +0.58%  "
$ws.Range("E34").Value = "  +1.23%  "
$ws.Range("D35").Value = "'2.52"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.63%  "
$ws.Range("D36").Value = "'17.49"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.92%  "
$ws.Range("E37").Value = "  -2.39%  "
$ws.Range("B38").Value = "Stellar"
$ws.Range("C38").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D38").Value = "'0.115"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.89%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "'0.103"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.99%  "
$ws.Range("E40").Value = "  -2.52%  "
$ws.Range("E41").Value = "  -3.36%  "
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("D43").Value = "1.974.26"
$ws.Range("E43").Value = "  +0.63%  "
$ws.Range("D44").Value = "'19.07"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.88%  "
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("E46").Value = "  -1.00%  "
$ws.Range("E47").Value = "  +2.27%  "
$ws.Range("D48").Value = "2.702.30"
$ws.Range("E48").Value = "  -0.90%  "
$ws.Range("D49").Value = "'97.21"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.86%  "
$ws.Range("D50").Value = "'66.99"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.80%  "
$ws.Range("D51").Value = "'72.92"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.81%  "
